$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 86

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "02/18/2026"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 9452.469999999999
$ws.Cells.Item($row, 3).Value = 0.2395311777718252
$ws.Cells.Item($row, 4).Value = 0.7604688222281748
$ws.Cells.Item($row, 5).Value = -314.14
$ws.Cells.Item($row, 6).Value = -35.5
$ws.Cells.Item($row, 7).Value = -23753.45
$ws.Cells.Item($row, 8).Value = -76.77
$ws.Cells.Item($row, 9).Value = -1138.69
$ws.Cells.Item($row, 10).Value = -33.46
$ws.Cells.Item($row, 11).Value = -24892.14
$ws.Cells.Item($row, 12).Value = -72.48
